$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 8255.93368014231
$ws.Range("F2").Value = 58.5854615201228

$ws.Range("C3").Value = 7310.57891135589
$ws.Range("F3").Value = 205.55779183404

$ws.Range("C4").Value = 4685.44230436666
$ws.Range("F4").Value = 73.8656730754888

$ws.Range("C5").Value = 4454.84946300448
$ws.Range("F5").Value = 58.6499785069484

$ws.Range("C6").Value = 7605.72970296246
$ws.Range("F6").Value = 231.799084748505

$ws.Range("C7").Value = 8790.20016646637
$ws.Range("F7").Value = 309.156014274735
